# Rename the "Glossary" tab to "Glossary&Definitions" and update the
# sheet's own "Section" column so its label matches the new tab name.
# Also moves the active-sheet/selection focus from "D&C - Procurement"
# over to the (now renamed) Glossary sheet, per "update script to
# include glossary page".

$wb = $excel.ActiveWorkbook

# --- D&C - Procurement: leave the last selection on E20, no longer the
#     active/focused tab ---
$wsProc = $wb.Worksheets.Item("D&C - Procurement")
$wsProc.Range("E20").Select() | Out-Null

# --- Glossary -> Glossary&Definitions ---
$wsGloss = $wb.Worksheets.Item("Glossary")
$wsGloss.Name = "Glossary&Definitions"

# The sheet's "Section" column (A2:A4) previously read "Glossary" for
# every row - update it to match the new sheet/tab name.
$wsGloss.Range("A2:A4").Value = "Glossary&Definitions"

# Make it the active tab/selection.
$wsGloss.Select() | Out-Null
$wsGloss.Range("D14").Select() | Out-Null

# Scroll the sheet-tab strip so the newly active tab is in view.
$win = $excel.ActiveWindow
$win.ScrollWorkbookTabs(1, 3) | Out-Null
